# Updated cryptos list on Wed Aug  9 18:39:56 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{Row=2; B='Bitcoin'; C='https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'; D='29.514.45'; E='  -1.10%  '; DText=$false},
    @{Row=3; B='Ethereum'; C='https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'; D='1.850.56'; E='  -0.73%  '; DText=$false},
    @{Row=4; B='TetherUSD'; C='https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'; D='0.9988'; E='  -0.08%  '; DText=$true},
    @{Row=5; B='BNB'; C='https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'; D='243.16'; E='  -1.63%  '; DText=$false},
    @{Row=6; B='XRP'; C='https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'; D='0.6538'; E='  +2.27%  '; DText=$true},
    @{Row=7; B='USDC'; C='https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'; D='0.9995'; E='  -0.05%  '; DText=$true},
    @{Row=8; B='Dogecoin'; C='https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'; D='0.07512'; E='  +0.08%  '; DText=$true},
    @{Row=9; B='Cardano'; C='https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'; D='0.2981'; E='  -0.91%  '; DText=$true},
    @{Row=10; B='Solana'; C='https://coinranking.com/coin/zNZHO_Sjf+solana-sol'; D='24.51'; E='  +1.04%  '; DText=$true},
    @{Row=11; B='TRON'; C='https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'; D='0.07634'; E='  -0.59%  '; DText=$true},
    @{Row=12; B='WrappedEther'; C='https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'; D='1.856.84'; E='  -0.86%  '; DText=$false},
    @{Row=13; B='Polkadot'; C='https://coinranking.com/coin/25W7FG7om+polkadot-dot'; D='5.023'; E='  -1.00%  '; DText=$true},
    @{Row=14; B='Polygon'; C='https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'; D='0.6853'; E='  -0.98%  '; DText=$true},
    @{Row=15; B='Litecoin'; C='https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'; D='83.74'; E='  -1.17%  '; DText=$true},
    @{Row=16; B='ShibaInu'; C='https://coinranking.com/coin/xz24e0BjL+shibainu-shib'; D='0.000009502'; E='  -0.18%  '; DText=$true},
    @{Row=17; B='Uniswap'; C='https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'; D='6.127'; E='  -0.15%  '; DText=$true},
    @{Row=18; B='WrappedBTC'; C='https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'; D='29.538.33'; E='  -0.89%  '; DText=$false},
    @{Row=19; B='WrappedliquidstakedEther2.0'; C='https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'; D='2.126.03'; E='  +1.61%  '; DText=$false},
    @{Row=20; B='BitcoinCash'; C='https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'; D='237.50'; E='  -1.42%  '; DText=$true},
    @{Row=21; B='Avalanche'; C='https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'; D='12.61'; E='  -0.80%  '; DText=$true},
    @{Row=22; B='Dai'; C='https://coinranking.com/coin/MoTuySvg7+dai-dai'; D='0.9993'; E='  -0.05%  '; DText=$true},
    @{Row=23; B='Chainlink'; C='https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'; D='7.700'; E='  +4.00%  '; DText=$true},
    @{Row=24; B='BinanceUSD'; C='https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'; D='1.000'; E='  -0.10%  '; DText=$true},
    @{Row=25; B='Monero'; C='https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'; D='157.03'; E='  -1.61%  '; DText=$false},
    @{Row=26; B='Stellar'; C='https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'; D='0.1422'; E='  -0.74%  '; DText=$true},
    @{Row=27; B='Cosmos'; C='https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'; D='8.505'; E='  -0.89%  '; DText=$true},
    @{Row=28; B='EthereumClassic'; C='https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'; D='17.82'; E='  -1.12%  '; DText=$false},
    @{Row=29; B='Hedera'; C='https://coinranking.com/coin/jad286TjB+hedera-hbar'; D='0.06041'; E='  +0.22%  '; DText=$true},
    @{Row=30; B='PancakeSwap'; C='https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'; D='1.487'; E='  -1.79%  '; DText=$true},
    @{Row=31; B='Toncoin'; C='https://coinranking.com/coin/67YlI0K1b+toncoin-ton'; D='1.254'; E='  -0.93%  '; DText=$false},
    @{Row=32; B='Filecoin'; C='https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; D='4.139'; E='  -0.33%  '; DText=$true},
    @{Row=33; B='InternetComputer(DFINITY)'; C='https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'; D='4.076'; E='  -1.77%  '; DText=$true},
    @{Row=34; B='ARBITRUM'; C='https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'; D='1.181'; E='  +1.65%  '; DText=$false},
    @{Row=35; B='LidoDAOToken'; C='https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'; D='1.855'; E='  -1.27%  '; DText=$false},
    @{Row=36; B='ImmutableX'; C='https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'; D='0.7248'; E='  -1.67%  '; DText=$false},
    @{Row=37; B='HuobiToken'; C='https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'; D='2.592'; E='  -0.89%  '; DText=$true},
    @{Row=38; B='MXToken'; C='https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'; D='2.802'; E='  -2.62%  '; DText=$false},
    @{Row=39; B='VeChain'; C='https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'; D='0.01784'; E='  -0.61%  '; DText=$true},
    @{Row=40; B='Maker'; C='https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'; D='1.202.10'; E='  -2.21%  '; DText=$false},
    @{Row=41; B='FraxShare'; C='https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'; D='6.248'; E='  -2.37%  '; DText=$true},
    @{Row=42; B='TrustWalletToken'; C='https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'; D='0.9092'; E='  -1.82%  '; DText=$true},
    @{Row=43; B='PaxDollar'; C='https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'; D='0.9989'; E='  -0.26%  '; DText=$true},
    @{Row=44; B='RocketPoolETH'; C='https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'; D='2.022.01'; E='  +0.88%  '; DText=$false},
    @{Row=45; B='Quant'; C='https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'; D='101.98'; E='  -0.63%  '; DText=$false},
    @{Row=46; B='Aave'; C='https://coinranking.com/coin/ixgUfzmLR+aave-aave'; D='66.49'; E='  -0.35%  '; DText=$true},
    @{Row=47; B='Aptos'; C='https://coinranking.com/coin/HGYj5JCv5+aptos-apt'; D='7.404'; E='  +10.10%  '; DText=$true},
    @{Row=48; B='BabyDogeCoin'; C='https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'; D='0.00000000121'; E='  -0.09%  '; DText=$true},
    @{Row=49; B='TheSandbox'; C='https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'; D='0.4056'; E='  -1.21%  '; DText=$true},
    @{Row=50; B='EnergySwap'; C='https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; D='9.089'; E='  -2.83%  '; DText=$true},
    @{Row=51; B='RenderToken'; C='https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'; D='1.664'; E='  +1.39%  '; DText=$true}
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $dCell = $ws.Cells.Item($r.Row, 4)
    if ($r.DText) {
        $dCell.NumberFormat = "@"
    }
    $dCell.Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
}
